# Added date fields to blank evaluation spreadsheet
# - Insert "EntryDate" column right after CompadrinoName (new column B)
# - Insert "ValidationDate" column right after ValidatedBy (new column D,
#   i.e. right after the EntryDate column that was just inserted)
# - Leave the selection on the newly-inserted ValidationDate column, as the
#   author would after inserting/inspecting it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data entry")

# --- Insert the two new blank columns, shifting existing data right -------
$ws.Columns("B:B").Insert()
$ws.Columns("D:D").Insert()

# --- Header labels for the new columns ------------------------------------
$ws.Range("B1").Value = "EntryDate"
$ws.Range("D1").Value = "ValidationDate"

# Match the bold header style used by the rest of row 1.
$ws.Range("B1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true

# --- New columns inherit the display width of their left neighbour --------
# (column B copies column A's 16-char width, column D copies column C's
# ~11.16-char width) and lose the old "best fit" flag, matching what Excel
# does when a column is inserted next to an explicitly-sized column.
$ws.Columns("B:B").ColumnWidth = 15.166666666666668
$ws.Columns("D:D").ColumnWidth = 10.25

# --- Leave the whole new ValidationDate column selected --------------------
$ws.Columns("D:D").Select()
